# "Generate Report for Archive"
#
# The localization-status report was regenerated:
#   * every cell whose status text was "Ready for handoff" is now
#     "In Translation" (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4)
#   * the two narrower "Status"-ish columns (Overview E:F, zh-cn C,
#     de-de C) shrink to fit the new, shorter status text
#
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- replace the status text -------------------------------------------------
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value     = "In Translation"
$dede.Range("C2:C4").Value     = "In Translation"

# --- shrink the columns that held that status text ---------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth     = 12.5   # column C (Status)
$dede.Columns.Item(3).ColumnWidth     = 12.5   # column C (Status)
